$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 21, pushing the existing rows 21-33 down to 22-34.
$ws.Rows("21:21").Insert()

# Populate the new row 21 with the weekly record (matches the pattern of
# the surrounding rows: same market/region/category/quality/unit/origin).
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44777
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112035
$ws.Range("G21").Value = "Bruselas (repollito)"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14400
$ws.Range("N21").Value = "$/malla 15 kilos"
$ws.Range("O21").Value = "Provincia de Quillota"
$ws.Range("P21").Value = 960
$ws.Range("Q21").Value = 15
$ws.Range("R21").Value = "Hortaliza"
